# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.006.62'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.235.14'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.91'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.18'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -6.76%  '
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -3.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.70'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -5.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0807'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.18'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -4.37%  '
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').Value = '2.577.22'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '2.236.77'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.820'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -3.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.50'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -4.67%  '
$ws.Range('D18').Value = '43.877.24'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').Value = '0.0₃0961'
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.09'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -9.09%  '
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.57'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('E24').Value = '  -5.88%  '
$ws.Range('E25').Value = '  -4.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.86'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -6.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.98'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.00'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.88'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.95'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -5.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0796'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -5.61%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.59'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.17'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.109'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -8.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.90'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -7.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.82'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -7.72%  '
$ws.Range('E41').Value = '  -8.75%  '
$ws.Range('E42').Value = '  -6.11%  '
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = '1.728.32'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '84.73'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('E46').Value = '  -4.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.59'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.91'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -4.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.04'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('E50').Value = '  -7.82%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '14.25'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +0.26%  '
